$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Cases")

# Update Run_Mode (column C) values for rows 3-6 (swap Y/N flags)
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "N"

# Move the active selection from H6 to K6 on the Test_Cases sheet
$ws.Range("K6").Select() | Out-Null
